$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

# Update the "address" column (D) values on the Child sheet (rows 2-21)
$ws.Range("D2").Value = "-5,6"
$ws.Range("D3").Value = "-9,-9"
$ws.Range("D4").Value = "-10,1"
$ws.Range("D5").Value = "7,2"
$ws.Range("D6").Value = "2,-2"
$ws.Range("D7").Value = "3,-9"
$ws.Range("D8").Value = "-9,-5"
$ws.Range("D9").Value = "-8,-7"
$ws.Range("D10").Value = "-10,-7"
$ws.Range("D11").Value = "-10,0"
$ws.Range("D12").Value = "7,-1"
$ws.Range("D13").Value = "-4,4"
$ws.Range("D14").Value = "8,0"
$ws.Range("D15").Value = "-7,-7"
$ws.Range("D16").Value = "0,-3"
$ws.Range("D17").Value = "-8,0"
$ws.Range("D18").Value = "-3,-7"
$ws.Range("D19").Value = "-2,-1"
$ws.Range("D20").Value = "-6,-9"
$ws.Range("D21").Value = "7,1"
